$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Metadata
$ws2 = $wb.Worksheets.Item(2)   # Concepts

# ---------------------------------------------------------------------------
# Sheet "Metadata"
# ---------------------------------------------------------------------------

# Version: 1.0.0 -> 0.1.0
$ws1.Range("B3").Value = "0.1.0"

# Status: active -> draft
$ws1.Range("B6").Value = "draft"

# Experimental: (blank) -> false   (must stay a text string, not boolean)
$ws1.Range("B7").Value = "'false"
$ws1.Range("B6").Copy()
$ws1.Range("B7").PasteSpecial(-4122)

# Date: 2025-11-28T01:24:36+00:00 -> 2025-12-26T14:13:58+00:00
$ws1.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description: (blank) -> Code system for categorizing types of nursing problems
$ws1.Range("B11").Value = "Code system for categorizing types of nursing problems"

# Case Sensitive: (blank) -> true   (must stay a text string, not boolean)
$ws1.Range("B14").Value = "'true"
$ws1.Range("B7").Copy()
$ws1.Range("B14").PasteSpecial(-4122)

# Count: 1 -> 3   (must stay a text string, not a number)
$ws1.Range("B21").Value = "'3"
$ws1.Range("B7").Copy()
$ws1.Range("B21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "Concepts"
# ---------------------------------------------------------------------------

# Fill in the Definition for the existing "nursing-diagnosis" row
$ws2.Range("D2").Value = "A clinical judgment about individual, family, or community responses to actual or potential health problems"
$ws2.Range("C2").Copy()
$ws2.Range("D2").PasteSpecial(-4122)

# New row 3: risk-diagnosis
$ws2.Range("A3").Value = "'1"
$ws2.Range("A2").Copy()
$ws2.Range("A3").PasteSpecial(-4122)

$ws2.Range("B3").Value = "risk-diagnosis"
$ws2.Range("B2").Copy()
$ws2.Range("B3").PasteSpecial(-4122)

$ws2.Range("C3").Value = "Risk Diagnosis"
$ws2.Range("C2").Copy()
$ws2.Range("C3").PasteSpecial(-4122)

$ws2.Range("D3").Value = "A clinical judgment about an individual's vulnerability to developing an undesirable health condition"
$ws2.Range("C2").Copy()
$ws2.Range("D3").PasteSpecial(-4122)

# New row 4: health-promotion
$ws2.Range("A4").Value = "'1"
$ws2.Range("A3").Copy()
$ws2.Range("A4").PasteSpecial(-4122)

$ws2.Range("B4").Value = "health-promotion"
$ws2.Range("B2").Copy()
$ws2.Range("B4").PasteSpecial(-4122)

$ws2.Range("C4").Value = "Health Promotion Diagnosis"
$ws2.Range("C2").Copy()
$ws2.Range("C4").PasteSpecial(-4122)

$ws2.Range("D4").Value = "A clinical judgment about motivation to increase wellbeing"
$ws2.Range("C2").Copy()
$ws2.Range("D4").PasteSpecial(-4122)
